# "Rework chapters on transformer" (#134) / Adresses #10
#
# The deck's two transformer-stack connector labels on slide 1 are
# relabelled from "N×" to "L×" (shapes 233/234, the TextBox callouts
# that sit on the "Nx"/"Lx" repeat-block connectors).
#
# The footer's cached "today" field (datetimeFigureOut) on the slide
# master and every slide layout also advances from 15/01/2023 to
# 28/01/2023 - PowerPoint recalculates that auto-update date field's
# cached text whenever the deck is subsequently saved, which is what
# happened when this edit was made a couple of weeks later.

$p = $ppt.ActivePresentation

# --- 1. Relabel the two "N×" connector callouts on slide 1 -----------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (($sh.Id -eq 233 -or $sh.Id -eq 234) -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq "N×") {
            $sh.TextFrame.TextRange.Text = "L×"
        }
    }
}

# --- 2. Refresh the cached footer date on the master + every layout --------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq "15/01/2023") {
                    $sh.TextFrame.TextRange.Text = "28/01/2023"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}
